$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.501.95'
$ws.Range('E2').Value = '  -4.16%  '
$ws.Range('D3').Value = '3.620.47'
$ws.Range('E3').Value = '  -4.59%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.95'
$ws.Range('E5').Value = '  -3.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.33'
$ws.Range('E6').Value = '  -1.48%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.611'
$ws.Range('E7').Value = '  -5.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.677'
$ws.Range('E9').Value = '  -8.31%  '
$ws.Range('E10').Value = '  -12.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.24'
$ws.Range('E11').Value = '  -6.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000253'
$ws.Range('E12').Value = '  -16.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.00'
$ws.Range('E13').Value = '  -7.67%  '
$ws.Range('D14').Value = '4.191.08'
$ws.Range('E14').Value = '  -4.37%  '
$ws.Range('D15').Value = '3.616.71'
$ws.Range('E15').Value = '  -4.42%  '
$ws.Range('E16').Value = '  -0.75%  '
$ws.Range('D17').Value = '67.230.97'
$ws.Range('E17').Value = '  -4.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.42'
$ws.Range('E18').Value = '  -6.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.28'
$ws.Range('E19').Value = '  -7.00%  '
$ws.Range('E20').Value = '  -7.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '395.53'
$ws.Range('E21').Value = '  -6.62%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.33'
$ws.Range('E22').Value = '  -8.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '85.23'
$ws.Range('E23').Value = '  -6.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.86'
$ws.Range('E24').Value = '  -8.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.34'
$ws.Range('E25').Value = '  -6.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.06'
$ws.Range('E26').Value = '  -0.64%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.37'
$ws.Range('E27').Value = '  -9.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.60'
$ws.Range('E28').Value = '  -12.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.01'
$ws.Range('E29').Value = '  -7.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.28'
$ws.Range('E30').Value = '  -6.95%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.79'
$ws.Range('E31').Value = '  -10.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '66.16'
$ws.Range('E32').Value = '  +1.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.94'
$ws.Range('E33').Value = '  -6.31%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.112'
$ws.Range('E34').Value = '  -7.24%  '
$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '588.71'
$ws.Range('E35').Value = '  -5.35%  '
$ws.Range('E36').Value = '  -7.23%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  -0.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.377'
$ws.Range('E39').Value = '  -8.95%  '
$ws.Range('D40').Value = '0.0₃0738'
$ws.Range('E40').Value = '  -20.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.134'
$ws.Range('E41').Value = '  -4.72%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.79'
$ws.Range('E42').Value = '  -10.82%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0413'
$ws.Range('E43').Value = '  -8.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.44'
$ws.Range('E44').Value = '  -13.95%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '26.67'
$ws.Range('E45').Value = '  +15.36%  '
$ws.Range('D46').Value = '2.688.65'
$ws.Range('E46').Value = '  -4.99%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.130'
$ws.Range('E47').Value = '  -5.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '139.79'
$ws.Range('E48').Value = '  -2.30%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.02'
$ws.Range('E49').Value = '  -8.11%  '
$ws.Range('B50').Value = 'WEMIXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.54'
$ws.Range('E50').Value = '  -8.10%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.42'
$ws.Range('E51').Value = '  -11.78%  '
